# Update cryptos list cell values per upstream diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '43.287.10'
$ws.Range('E2').Value = '  +2.63%  '
$ws.Range('D3').Value = '2.302.35'
$ws.Range('E3').Value = '  +1.53%  '
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '310.72'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.46%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '102.40'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +5.30%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.532'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +1.08%  '
$ws.Range('E8').Value = '  +0.01%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.528'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +7.68%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '35.55'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +1.09%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0810'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +2.63%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '6.95'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +1.10%  '
$ws.Range('D14').Value = '2.661.51'
$ws.Range('E14').Value = '  +1.59%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '14.96'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +1.79%  '
$ws.Range('D16').Value = '2.294.53'
$ws.Range('E16').Value = '  +1.61%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.808'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +2.07%  '
$ws.Range('D18').Value = '43.199.84'
$ws.Range('E18').Value = '  +2.63%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '12.22'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -0.49%  '
$ws.Range('D20').Value = '0.0₃0930'
$ws.Range('E20').Value = '  +2.75%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.16'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +2.54%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '67.96'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.23%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '240.80'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +1.54%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.61'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +1.14%  '
$ws.Range('E25').Value = '  +1.70%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.999'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -0.07%  '
$ws.Range('E27').Value = '  -1.88%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '24.76'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +5.26%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '36.61'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -2.50%  '
$ws.Range('E30').Value = '  +3.30%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '9.61'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +0.27%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '169.52'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +4.34%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '5.25'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +0.04%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.54'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +6.85%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '17.73'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +0.54%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.0739'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +0.51%  '
$ws.Range('E38').Value = '  -2.89%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.88'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +2.61%  '
$ws.Range('E40').Value = '  +1.80%  '
$ws.Range('E41').Value = '  +0.92%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '4.34'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +6.17%  '
$ws.Range('E43').Value = '  -1.38%  '
$ws.Range('B44').Value = 'VeChain'
$ws.Range('C44').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0288'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +2.57%  '
$ws.Range('B45').Value = 'Maker'
$ws.Range('C45').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D45').Value = '1.962.04'
$ws.Range('E45').Value = '  +0.60%  '
$ws.Range('B46').Value = 'EnergySwap'
$ws.Range('C46').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '19.14'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +0.95%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.98'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +2.13%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '9.86'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -1.10%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '55.21'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +2.49%  '
$ws.Range('E50').Value = '  +0.98%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.57'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +6.92%  '
